# Commit: "adding 4th edition edited list of tumor"
#
# 1. Rename sheet "Final_List" -> "Edited_Tumors"
# 2. Make "Edited_Tumors" the active/selected sheet (moves tabSelected +
#    sets workbook-level activeTab)
# 3. Consolidate the three per-column ("B", "C", "D") "contains Yes" red-fill
#    conditional format rules on the "Generated" sheet into a single rule
#    applied across B:D, keeping the rule/format that used to live on column D
#    (dxfId 0 / priority 1) and widening it to cover columns B:D.

$wb = $excel.ActiveWorkbook

# --- Rename the sheet ------------------------------------------------------
$editedSheet = $wb.Worksheets.Item("Final_List")
$editedSheet.Name = "Edited_Tumors"

# --- Consolidate conditional formatting on "Generated" ---------------------
$gen = $wb.Worksheets.Item("Generated")

# Drop the column-B and column-C rules entirely.
$gen.Range("B1:B1048576").FormatConditions.Delete()
$gen.Range("C1:C1048576").FormatConditions.Delete()

# Re-scope the remaining (column-D) rule so it covers B:D instead of just D,
# keeping its existing dxf/fill.
$keepRule = $gen.Range("D1:D1048576").FormatConditions.Item(1)
$keepRule.ModifyAppliesToRange($gen.Range("B1:D1048576"))
$keepRule.Formula1 = 'NOT(ISERROR(SEARCH("Yes",B1)))'

# --- Activate the renamed sheet (becomes the selected tab) -----------------
$editedSheet.Activate()
